# Apply the authored changes to TxRecord.xlsx:
#  - DBS sheet (sheet2): append two new lookup rows (findByCalDate / findByCustNo)
#  - Update remembered selections on both sheets
#  - Keep DBS as the active/selected tab

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("DBD")
$ws2 = $wb.Worksheets.Item("DBS")

# --- DBS: add the two new rows -------------------------------------------
# Row 10
$ws2.Range("A11").Value = "findByCustNo"
$ws2.Range("B10").Value = "CalDate >= ,AND CalDate <= ,AND BrNo = ,AND ImportFg ="
$ws2.Range("B11").Value = "CalDate >= ,AND CalDate <= ,AND BrNo = ,AND LockCustNo = ,AND ImportFg ="
$ws2.Range("A10").Value = "findByCalDate"
$ws2.Range("C10").Value = "CreateDate asc"
$ws2.Range("C11").Value = "CreateDate asc"

# --- Update remembered selection / scroll position for DBD ----------------
$ws1.Activate()
$ws1.Range("C32").Select()

# --- Reactivate DBS (it must remain the active/selected tab) and update its
#     remembered selection ---------------------------------------------------
$ws2.Activate()
$ws2.Range("B14").Select()
